$wb = $excel.ActiveWorkbook

# Target tab order: 总计, 2021-Q4, 2021-Q3
# Step 1: move "总计" to the very front (before "2021-Q3")
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ3 = $wb.Worksheets.Item("2021-Q3")
$wsTotal.Move($wsQ3)

# Step 2: move "2021-Q3" to the very end (after "2021-Q4")
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ3 = $wb.Worksheets.Item("2021-Q3")
$wsQ3.Move($null, $wsQ4)

# "2021-Q3" keeps the selected/active tab (it was the tabSelected sheet before the edit)
$wsQ3 = $wb.Worksheets.Item("2021-Q3")
$wsQ3.Activate()
